# Apply the "Add CostTypeExportWizard" rework to ElementDefinitionTest.xlsx
#
# Summary of functional changes:
#  - Sheet "InterfaceTypes"  -> renamed "CostTypes"
#  - Sheet "InterfaceEnds"   -> renamed "CostEquipments"
#  - Sheet "Interfaces"      -> deleted entirely
#  - Header sheet: Structural Element Name value "BATTERY1" -> "BATTERY"
#  - CostTypes sheet header label updated to "VirSat IO Sheet for Cost Types"
#    and column header "Interface Type Name" -> "Cost Type Name"
#  - CostEquipments sheet header label updated to "VirSat IO Sheet for CostEquipments",
#    column headers "InterfaceEnd Name" -> "CostEquipment Name", "Interface Type" -> "Cost Type"
#  - CostEquipments data column (previously "Interface Type" values) updated:
#       HILL -> BIII
#       BILL -> BAAA
#       KILL -> BUUU (used by both POW_OUT rows and the trailing D9 cell)

$wb = $excel.ActiveWorkbook

# --- Delete the "Interfaces" worksheet ---------------------------------
$wsInterfaces = $wb.Worksheets.Item("Interfaces")
$wsInterfaces.Delete()

# --- Rename the remaining data sheets -----------------------------------
$wsTypes = $wb.Worksheets.Item("InterfaceTypes")
$wsTypes.Name = "CostTypes"

$wsEnds = $wb.Worksheets.Item("InterfaceEnds")
$wsEnds.Name = "CostEquipments"

# --- Header sheet: update the Structural Element Name value ------------
$wsHeader = $wb.Worksheets.Item("Header")
$wsHeader.Range("B6").Value = "BATTERY"

# --- CostTypes sheet: update titles/labels ------------------------------
$wsTypes.Range("A2").Value = "VirSat IO Sheet for Cost Types"
$wsTypes.Range("C4").Value = "Cost Type Name"

# --- CostEquipments sheet: update titles/labels and data ----------------
$wsEnds.Range("A2").Value = "VirSat IO Sheet for CostEquipments"
$wsEnds.Range("C4").Value = "CostEquipment Name"
$wsEnds.Range("D4").Value = "Cost Type"

$wsEnds.Range("D5").Value = "BIII"
$wsEnds.Range("D6").Value = "BAAA"
$wsEnds.Range("D7").Value = "BUUU"
$wsEnds.Range("D8").Value = "BUUU"
$wsEnds.Range("D9").Value = "BUUU"
